$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 6 data: function name, filename, and line number
$ws.Range("A6").Value = "isDeviceOnline"
$ws.Range("C6").Value = "/home/rdkv-core/cov/cov-analysis-linux64-2023.6.0/bin/device/entservices-softwareupdate/MaintenanceManager/MaintenanceManager.cpp"
$ws.Range("D6").Value = 1323

# Update the active cell / selection shown when the sheet was saved
$ws.Range("E9").Select()
